# Fixed (see issue #93)
# Rename the transformer-type template column headers on the
# "transformer_types" sheet to the shorter attribute names used by GridCal.

$wb = $excel.ActiveWorkbook

$wsBranch = $wb.Worksheets.Item("branch")
$wsXf     = $wb.Worksheets.Item("transformer_types")

$wsXf.Range("C1").Value = "HV"
$wsXf.Range("D1").Value = "LV"
$wsXf.Range("E1").Value = "rating"
$wsXf.Range("F1").Value = "Pcu"
$wsXf.Range("G1").Value = "Pfe"
$wsXf.Range("H1").Value = "I0"
$wsXf.Range("I1").Value = "Vsc"

# Restore the selection on the "branch" sheet before switching away from it.
[void]$wsBranch.Range("G5").Select()

# Leave the "transformer_types" sheet active/selected, matching the saved
# workbook view.
$wsXf.Activate()
[void]$wsXf.Range("D5").Select()
